$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row 39 (pushing the old row 39 "[TABLE_END]" down to row 40),
# carrying formatting down from row 38 (left/center vertical align, wrap).
# ---------------------------------------------------------------------------
$ws.Rows.Item(39).Insert()

# Row height for the new changelog row
$ws.Rows.Item(39).RowHeight = 108.75

# Column A: text id key
$ws.Cells.Item(39, 1).Value = "update_info"

# Column B: Chinese changelog text (plain, no rich-text runs)
$cn = "最近更新`n版本号：13.0.0`n发布日期：2025-04-01`n核心更新：`n多语言系统增强`n• 新增日语语言支持（801bda8）`n• 扩展基础语言框架（7ec452e）"
$ws.Cells.Item(39, 2).Value = $cn

# Column C: English changelog text (plain, no rich-text runs)
$en = "Recent Update`nVersion Number: 13.0.0`nRelease Date: April 1, 2025`nCore Updates:`nEnhancement of the Multilingual System`n• Added Japanese language support (801bda8)`n• Extended the basic language framework (7ec452e)"
$ws.Cells.Item(39, 3).Value = $en

# Column D: Japanese changelog text with mixed-font rich-text runs
$run1 = "最近のアップデート`nバージョン番号：13.0.0`nリリース日：2025 年 4 月 1 日`n核心アップデート：`n多言語システムの強化`n"
$run2 = "・日本語の言語サポートを新たに追加しました"
$run3 = "（801bda8）`n"
$run4 = "・基本的な言語フレームワークを拡張しました"
$run5 = "（7ec452e）"
$jp = $run1 + $run2 + $run3 + $run4 + $run5

$dCell = $ws.Cells.Item(39, 4)
$dCell.Value = $jp

$pos = 1 + $run1.Length
$chars2 = $dCell.Characters($pos, $run2.Length)
$chars2.Font.Name = "Yu Gothic"
$chars2.Font.Size = 11

$pos = $pos + $run2.Length
$chars3 = $dCell.Characters($pos, $run3.Length)
$chars3.Font.Name = "等线"
$chars3.Font.Size = 11

$pos = $pos + $run3.Length
$chars4 = $dCell.Characters($pos, $run4.Length)
$chars4.Font.Name = "Yu Gothic"
$chars4.Font.Size = 11

$pos = $pos + $run4.Length
$chars5 = $dCell.Characters($pos, $run5.Length)
$chars5.Font.Name = "等线"
$chars5.Font.Size = 11

# ---------------------------------------------------------------------------
# Sheet view: scroll / selection around the newly added row.
# ---------------------------------------------------------------------------
$ws.Range("B39").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "edit applied"
